$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.837.54'
$ws.Range("E2").Value = '  -1.16%  '
$ws.Range("D3").Value = '1.872.50'
$ws.Range("E3").Value = '  -1.37%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.52'
$ws.Range("E5").Value = '  -1.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5373'
$ws.Range("E7").Value = '  +2.64%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3758'
$ws.Range("E8").Value = '  -1.26%  '
$ws.Range("E9").Value = '  -1.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.58'
$ws.Range("E10").Value = '  +0.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8876'
$ws.Range("E11").Value = '  -1.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08143'
$ws.Range("E12").Value = '  -0.22%  '
$ws.Range("D13").Value = '1.853.09'
$ws.Range("E13").Value = '  +1.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.36'
$ws.Range("E14").Value = '  -2.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.270'
$ws.Range("E15").Value = '  -1.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.74'
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008544'
$ws.Range("E18").Value = '  -1.24%  '
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("D20").Value = '26.858.50'
$ws.Range("E20").Value = '  -1.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.976'
$ws.Range("E21").Value = '  -2.36%  '
$ws.Range("E22").Value = '  -0.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.384'
$ws.Range("E23").Value = '  -0.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '147.20'
$ws.Range("E24").Value = '  -1.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.263'
$ws.Range("E25").Value = '  -2.67%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.736'
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("E27").Value = '  -1.00%  '
$ws.Range("E28").Value = '  -1.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.733'
$ws.Range("E29").Value = '  -1.86%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.596'
$ws.Range("E30").Value = '  -5.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09141'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8108'
$ws.Range("E32").Value = '  +2.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04972'
$ws.Range("E33").Value = '  -1.50%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.171'
$ws.Range("E34").Value = '  -4.14%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.982'
$ws.Range("E35").Value = '  +0.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6104'
$ws.Range("E36").Value = '  +7.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.600'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.178'
$ws.Range("E38").Value = '  -5.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01957'
$ws.Range("E39").Value = '  -1.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.074'
$ws.Range("E40").Value = '  -0.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.584'
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.834'
$ws.Range("E42").Value = '  -1.90%  '
$ws.Range("B43").Value = 'Decentraland'
$ws.Range("C43").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5189'
$ws.Range("E43").Value = '  +6.28%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '115.87'
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1498'
$ws.Range("E45").Value = '  -0.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9999'
$ws.Range("E46").Value = '  -0.28%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.965'
$ws.Range("E47").Value = '  -1.71%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.633'
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '37.64'
$ws.Range("E49").Value = '  -2.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06029'
$ws.Range("E50").Value = '  +1.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '62.05'
$ws.Range("E51").Value = '  -2.82%  '
